# Refresh the crypto price/volume snapshot (cryptos.xlsx) to match the
# GitHub Actions scrape from Sun Oct 8 08:12:12 UTC 2023.
#
# All Coin/Link/Price/Volume(1h) cells are plain text (inline strings) in
# the source workbook, so every write below targets text content. For
# "Price" cells whose new value happens to look like a clean decimal
# number (e.g. "211.78"), a leading apostrophe forces Excel to keep the
# cell as text instead of silently re-typing it as a Number (which would
# also strip meaningful trailing zeros, e.g. "230.50" -> 230.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.906.12"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.628.96"
$ws.Range("E3").Value = "  -0.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'211.78"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.04%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'23.31"
$ws.Range("E8").Value = "  -0.83%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.68%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.19%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0881"
$ws.Range("E11").Value = "  +0.69%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.861.41"
$ws.Range("E12").Value = "  -0.46%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.627.34"
$ws.Range("E13").Value = "  -0.50%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -1.59%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.561"
$ws.Range("E15").Value = "  -2.24%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.68"
$ws.Range("E16").Value = "  -0.35%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.906.77"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'230.50"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19 - was ShibaInu, is now Chainlink (rows 19/20 swapped order)
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.67"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20 - was Chainlink, is now ShibaInu
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.02%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.34"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "'10.23"
$ws.Range("E23").Value = "  -4.89%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.29%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'154.76"
$ws.Range("E25").Value = "  +2.04%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.10%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -0.90%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.05%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.87%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.46%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.33%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.402.12"
$ws.Range("E34").Value = "  +0.28%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.34%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +9.99%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.14%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +2.27%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +0.46%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "'0.862"
$ws.Range("E40").Value = "  -2.79%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -0.19%  "

# Row 43 - RenderToken
$ws.Range("D43").Value = "'1.84"
$ws.Range("E43").Value = "  +0.75%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'66.20"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45 - FraxShare
$ws.Range("D45").Value = "'5.49"
$ws.Range("E45").Value = "  +0.98%  "

# Row 46 - MXToken
$ws.Range("D46").Value = "'2.19"
$ws.Range("E46").Value = "  -0.30%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.770.66"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'88.06"
$ws.Range("E48").Value = "  +0.21%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -2.23%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +0.87%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.36%  "
